# Adds season-record columns (Wins / Losses / Ties) to the player table.
# Mirrors the upstream fix described in the commit message: the scraper
# now also pulls the team's season record and appends it as three new
# columns (AD, AE, AF) alongside the existing per-player stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - same bold/centered/bordered formatting as the
# other header cells (style index 1 in the original sheet).
$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Season record is the same for every player row (2-62): 76 wins,
# 86 losses, 0 ties.
$wins = 76
$losses = 86
$ties = 0

for ($row = 2; $row -le 62; $row++) {
    $ws.Cells.Item($row, 30).Value2 = $wins    # AD
    $ws.Cells.Item($row, 31).Value2 = $losses  # AE
    $ws.Cells.Item($row, 32).Value2 = $ties    # AF
}
